# Correção da impressão das cores nas curvas de luz
# As cores do espectro estavam invertidas
# Updates the D_lambda [ppm] values in column C (rows 2-59) of Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @(
    12359.7663786581,
    12365.52392646739,
    12386.16764613731,
    12383.78269112561,
    12393.38329307316,
    12403.8113467031,
    12361.45276949163,
    12334.56858929871,
    12375.39338907934,
    12336.28930582964,
    12360.79728425732,
    12356.07044203646,
    12355.20782869404,
    12363.90302125456,
    12352.23787712647,
    12344.47452774534,
    12332.88999475801,
    12334.420208576,
    12333.99086895182,
    12329.76728460067,
    12325.52286418853,
    12317.2638704897,
    12311.71614746707,
    12308.63988667974,
    12296.9871061458,
    12289.18074932917,
    12283.00161538054,
    12277.49691427071,
    12273.07547475431,
    12268.20311668031,
    12262.89140110104,
    12260.62749882806,
    12253.44332547795,
    12250.12565907346,
    12244.57685625546,
    12240.46338575469,
    12240.46227324216,
    12238.22290210497,
    12231.43198938337,
    12230.22580934763,
    12228.72067064173,
    12224.46111875475,
    12222.49290238164,
    12222.03529727284,
    12220.16294890738,
    12220.9601840434,
    12217.90965531133,
    12213.27088713808,
    12210.50225598752,
    12205.879261733,
    12203.9038825722,
    12199.27704047385,
    12193.15113390462,
    12186.92963477264,
    12180.7046047302,
    12182.21530017216,
    12175.12257087294,
    12169.82407055778
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $newValues[$i]
}
